$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Range, $Value)
    $Range.NumberFormat = "@"
    $Range.Value = $Value
    $Range.ClearFormats()
}

Set-TextValue $ws.Range("D2") '60.391.87'
Set-TextValue $ws.Range("E2") '  +1.91%  '
Set-TextValue $ws.Range("D3") '2.585.27'
Set-TextValue $ws.Range("E3") '  +1.91%  '
Set-TextValue $ws.Range("E4") '  -0.40%  '
Set-TextValue $ws.Range("D5") '506.40'
Set-TextValue $ws.Range("E5") '  +0.29%  '
Set-TextValue $ws.Range("D6") '154.02'
Set-TextValue $ws.Range("E6") '  -1.35%  '
Set-TextValue $ws.Range("D7") '0.999'
Set-TextValue $ws.Range("E7") '  +0.60%  '
Set-TextValue $ws.Range("E8") '  -8.12%  '
Set-TextValue $ws.Range("D9") '2.591.66'
Set-TextValue $ws.Range("E9") '  +0.34%  '
Set-TextValue $ws.Range("D10") '6.57'
Set-TextValue $ws.Range("E10") '  +6.46%  '
Set-TextValue $ws.Range("E11") '  +0.77%  '
Set-TextValue $ws.Range("D12") '0.347'
Set-TextValue $ws.Range("E12") '  +1.85%  '
Set-TextValue $ws.Range("E13") '  +0.84%  '
Set-TextValue $ws.Range("D14") '3.036.66'
Set-TextValue $ws.Range("E14") '  +0.74%  '
Set-TextValue $ws.Range("D15") '60.400.83'
Set-TextValue $ws.Range("E15") '  +2.13%  '
Set-TextValue $ws.Range("E16") '  -1.23%  '
Set-TextValue $ws.Range("E17") '  +1.72%  '
Set-TextValue $ws.Range("D18") '2.586.51'
Set-TextValue $ws.Range("E18") '  +0.40%  '
Set-TextValue $ws.Range("D19") '4.80'
Set-TextValue $ws.Range("E19") '  +0.85%  '
Set-TextValue $ws.Range("D20") '345.44'
Set-TextValue $ws.Range("E20") '  +2.87%  '
Set-TextValue $ws.Range("E21") '  +0.87%  '
Set-TextValue $ws.Range("E22") '  +1.54%  '
Set-TextValue $ws.Range("D23") '0.997'
Set-TextValue $ws.Range("E23") '  -0.92%  '
Set-TextValue $ws.Range("D24") '59.98'
Set-TextValue $ws.Range("E24") '  +0.28%  '
Set-TextValue $ws.Range("E25") '  +1.28%  '
Set-TextValue $ws.Range("E26") '  -0.23%  '
Set-TextValue $ws.Range("D27") '0.997'
Set-TextValue $ws.Range("E27") '  -0.32%  '
Set-TextValue $ws.Range("D28") '0.0₃0845'
Set-TextValue $ws.Range("E28") '  +1.99%  '
Set-TextValue $ws.Range("D29") '7.32'
Set-TextValue $ws.Range("E29") '  -1.10%  '
Set-TextValue $ws.Range("E30") '  +0.33%  '
Set-TextValue $ws.Range("D31") '19.33'
Set-TextValue $ws.Range("E31") '  +0.22%  '
Set-TextValue $ws.Range("D32") '153.56'
Set-TextValue $ws.Range("E32") '  -2.10%  '
Set-TextValue $ws.Range("E33") '  -1.04%  '
Set-TextValue $ws.Range("E34") '  +3.57%  '
Set-TextValue $ws.Range("E35") '  +1.79%  '
Set-TextValue $ws.Range("E36") '  -0.06%  '
Set-TextValue $ws.Range("D37") '0.858'
Set-TextValue $ws.Range("E37") '  +12.79%  '
Set-TextValue $ws.Range("D38") '0.847'
Set-TextValue $ws.Range("E38") '  -0.44%  '
Set-TextValue $ws.Range("E39") '  +1.42%  '
Set-TextValue $ws.Range("E40") '  +0.35%  '
Set-TextValue $ws.Range("D41") '35.89'
Set-TextValue $ws.Range("E41") '  +2.02%  '
Set-TextValue $ws.Range("D42") '295.89'
Set-TextValue $ws.Range("E42") '  +1.39%  '
Set-TextValue $ws.Range("E43") '  -1.93%  '
Set-TextValue $ws.Range("D45") '0.0557'
Set-TextValue $ws.Range("E45") '  -0.80%  '
Set-TextValue $ws.Range("D46") '0.997'
Set-TextValue $ws.Range("E46") '  +0.73%  '
Set-TextValue $ws.Range("D47") '19.81'
Set-TextValue $ws.Range("E47") '  +3.13%  '
Set-TextValue $ws.Range("D48") '4.84'
Set-TextValue $ws.Range("E48") '  -0.03%  '
Set-TextValue $ws.Range("E49") '  -0.97%  '
Set-TextValue $ws.Range("D50") '10.31'
Set-TextValue $ws.Range("E50") '  +0.57%  '
Set-TextValue $ws.Range("D51") '2.000.32'
Set-TextValue $ws.Range("E51") '  -0.24%  '
